$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A3").Value = 301.72000000000003
$wsSummary.Range("E3").Value = 301.72000000000003

# --- Sheet: Repayment schedule -----------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Row 11 (installment #9)
$wsRepay.Range("B11").Value = 14
$wsRepay.Range("C11").Value = 42157
$wsRepay.Range("F11").Value = 842.84
$wsRepay.Range("G11").Value = 2551.96
$wsRepay.Range("H11").Value = 15.63

# Row 12 (installment #10)
$wsRepay.Range("B12").Value = 14
$wsRepay.Range("F12").Value = 846.72
$wsRepay.Range("G12").Value = 1705.24
$wsRepay.Range("H12").Value = 11.75

# Row 13 (installment #11)
$wsRepay.Range("F13").Value = 850.62
$wsRepay.Range("G13").Value = 854.62
$wsRepay.Range("H13").Value = 7.85

# Row 14 (installment #12)
$wsRepay.Range("F14").Value = 854.62
$wsRepay.Range("H14").Value = 3.93
$wsRepay.Range("K14").Value = 858.55
$wsRepay.Range("Q14").Value = 858.55

# --- Selections / active sheet ------------------------------------------
# Select on the non-final-active sheets first, then finish on "Summary"
# so it ends up as the active tab (activeTab=1, tabSelected on Summary).

$wsToClient2 = $wb.Worksheets.Item("ToClient2")
$wsToClient2.Range("B4").Select()

$wsRepay.Range("M6").Select()

$wsSummary.Range("C2:D2").Select()
